# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (with per-fund holdings detail) right
# after "总计" and before "2022-Q3", and prepends a 2022-Q4 summary row
# to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)   # 总计
$q3    = $wb.Worksheets.Item(2)   # 2022-Q3 (current 2nd sheet; new sheet goes before it)

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet positioned before "2022-Q3"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Copy the header/index cell formatting from the "总计" sheet (style
# index already used for bold+bordered header / index cells) onto the
# new sheet so new cells share the same look instead of default style.
$total.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Copy()
$q4.Range("A2:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Fund holdings detail for 2022-Q4. Columns D/E/F/G are stored as text
# (matching the other quarter sheets), A/H are numeric.
$q4data = @(
  ,@("011486","博时创新精选混合A","9.23","92.26","3.58","0.3304",5)
  ,@("012082","博时数字经济混合A","4.95","91.76","4.89","0.2421",2)
  ,@("050022","博时回报混合","3.91","76.56","3.89","0.1521",5)
  ,@("014036","博时成长回报混合A","3.64","92.08","3.38","0.1230",9)
  ,@("673020","西部利得成长精选灵活配置混合","1.16","80.98","6.20","0.0719",3)
  ,@("011487","博时创新精选混合C","1.01","92.26","3.58","0.0362",5)
  ,@("014600","博时回报严选混合A","0.89","93.75","3.56","0.0317",7)
  ,@("014037","博时成长回报混合C","0.78","92.08","3.38","0.0264",9)
  ,@("015880","中欧小盘成长混合A","2.54","92.29","0.99","0.0251",6)
  ,@("012083","博时数字经济混合C","0.28","91.76","4.89","0.0137",2)
  ,@("015881","中欧小盘成长混合C","1.31","92.29","0.99","0.0130",6)
  ,@("014601","博时回报严选混合C","0.05","93.75","3.56","0.0018",7)
)

for ($i = 0; $i -lt $q4data.Length; $i++) {
  $row = $q4data[$i]
  $r = $i + 2

  $q4.Cells.Item($r, 1).Value = $i

  $q4.Cells.Item($r, 2).NumberFormat = "@"
  $q4.Cells.Item($r, 2).Value = $row[0]

  $q4.Cells.Item($r, 3).NumberFormat = "@"
  $q4.Cells.Item($r, 3).Value = $row[1]

  $q4.Cells.Item($r, 4).NumberFormat = "@"
  $q4.Cells.Item($r, 4).Value = $row[2]

  $q4.Cells.Item($r, 5).NumberFormat = "@"
  $q4.Cells.Item($r, 5).Value = $row[3]

  $q4.Cells.Item($r, 6).NumberFormat = "@"
  $q4.Cells.Item($r, 6).Value = $row[4]

  $q4.Cells.Item($r, 7).NumberFormat = "@"
  $q4.Cells.Item($r, 7).Value = $row[5]

  $q4.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: push existing rows down one and insert a
#    new 2022-Q4 summary row at the top of the data.
# ---------------------------------------------------------------------

# Give the about-to-exist row 4 the same formatting as row 3 (row 4
# doesn't exist yet in the sheet, so it has no style of its own).
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift 2022-Q2 (currently row 2) down to row 4, and 2022-Q3 (currently
# row... wait, already row 2 as well originally) -- do this via plain
# value rewrites from the bottom up so we don't clobber data we still
# need to read.
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0.19

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 7
$total.Cells.Item(3,4).Value = 0.65

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 12
$total.Cells.Item(2,4).Value = 1.07

# Keep the workbook's active tab as it was before this edit ("总计").
$total.Activate()
